$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(8, 8).Value = 80077.60000000001
$ws.Cells.Item(8, 9).Value = 85797.28999999999
$ws.Cells.Item(8, 10).Value = 2
$ws.Cells.Item(8, 11).Value = 257391.87
$ws.Cells.Item(8, 12).Value = 6
$ws.Cells.Item(8, 13).Value = -257252.87
$ws.Cells.Item(8, 14).Value = -284
$ws.Cells.Item(9, 8).Value = 689.1539
$ws.Cells.Item(9, 9).Value = 737.5454999999999
$ws.Cells.Item(9, 10).Value = 423
$ws.Cells.Item(9, 11).Value = 737.5454999999999
$ws.Cells.Item(9, 12).Value = 423
$ws.Cells.Item(9, 13).Value = -568.5454999999999
$ws.Cells.Item(9, 14).Value = -761
$ws.Cells.Item(18, 8).Value = 2777.3333
$ws.Cells.Item(18, 10).Value = 4982.5
$ws.Cells.Item(18, 12).Value = 4982.5
$ws.Cells.Item(18, 14).Value = -5550.5
$ws.Cells.Item(33, 8).Value = 457.81818
$ws.Cells.Item(33, 9).Value = 456
$ws.Cells.Item(33, 11).Value = 456
$ws.Cells.Item(33, 13).Value = -227
$ws.Cells.Item(43, 8).Value = 6111
$ws.Cells.Item(43, 10).Value = 5638.5
$ws.Cells.Item(43, 12).Value = 5638.5
$ws.Cells.Item(43, 14).Value = -5776.5
$ws.Cells.Item(64, 8).Value = 5375.625
$ws.Cells.Item(64, 9).Value = 3799.8
$ws.Cells.Item(64, 11).Value = 3799.8
$ws.Cells.Item(64, 13).Value = -3551.8
$ws.Cells.Item(67, 8).Value = 5375.625
$ws.Cells.Item(67, 9).Value = 3799.8
$ws.Cells.Item(67, 11).Value = 3799.8
$ws.Cells.Item(67, 13).Value = -2941.8
$ws.Cells.Item(74, 8).Value = 4820
$ws.Cells.Item(74, 9).Value = 3937.1428
$ws.Cells.Item(74, 11).Value = 3937.1428
$ws.Cells.Item(74, 13).Value = -3001.1428
$ws.Cells.Item(77, 8).Value = 4820
$ws.Cells.Item(77, 9).Value = 3937.1428
$ws.Cells.Item(77, 11).Value = 19685.714
$ws.Cells.Item(77, 13).Value = -15005.714
$ws.Cells.Item(80, 8).Value = 1520.5883
$ws.Cells.Item(80, 9).Value = 1754.4286
$ws.Cells.Item(80, 11).Value = 5263.2858
$ws.Cells.Item(80, 13).Value = -4265.2858
$ws.Cells.Item(83, 8).Value = 1520.5883
$ws.Cells.Item(83, 9).Value = 1754.4286
$ws.Cells.Item(83, 11).Value = 15789.8574
$ws.Cells.Item(83, 13).Value = -10797.8574
$ws.Cells.Item(96, 8).Value = 1603.6923
$ws.Cells.Item(96, 9).Value = 1703.375
$ws.Cells.Item(96, 10).Value = 1444.2
$ws.Cells.Item(96, 11).Value = 5110.125
$ws.Cells.Item(96, 12).Value = 4332.6
$ws.Cells.Item(96, 13).Value = -3737.125
$ws.Cells.Item(96, 14).Value = -7078.6
$ws.Cells.Item(100, 8).Value = 2132.111
$ws.Cells.Item(100, 9).Value = 1798.5
$ws.Cells.Item(100, 10).Value = 2799.3333
$ws.Cells.Item(100, 11).Value = 1798.5
$ws.Cells.Item(100, 12).Value = 2799.3333
$ws.Cells.Item(100, 13).Value = -1257.5
$ws.Cells.Item(100, 14).Value = -3881.3333
$ws.Cells.Item(113, 8).Value = 8396.666999999999
$ws.Cells.Item(113, 10).Value = 9482.5
$ws.Cells.Item(113, 12).Value = 9482.5
$ws.Cells.Item(113, 14).Value = -15990.5
$ws.Cells.Item(132, 8).Value = 3564.7727
$ws.Cells.Item(132, 9).Value = 3496.4285
$ws.Cells.Item(132, 11).Value = 10489.2855
$ws.Cells.Item(132, 13).Value = -7959.2855
$ws.Cells.Item(137, 8).Value = 1883.8889
$ws.Cells.Item(137, 10).Value = 1577.8334
$ws.Cells.Item(137, 12).Value = 4733.5002
$ws.Cells.Item(137, 14).Value = -9833.5002

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(19, 8).Value = 4000
$ws.Cells.Item(19, 9).Value = 4000
$ws.Cells.Item(19, 11).Value = 4000
$ws.Cells.Item(19, 13).Value = -3771
$ws.Cells.Item(32, 8).Value = 1386.4783
$ws.Cells.Item(32, 9).Value = 1280.4286
$ws.Cells.Item(32, 11).Value = 1280.4286
$ws.Cells.Item(32, 13).Value = -993.4286
$ws.Cells.Item(74, 8).Value = 2421.5715
$ws.Cells.Item(74, 9).Value = 1738.75
$ws.Cells.Item(74, 11).Value = 1738.75
$ws.Cells.Item(74, 13).Value = -864.75
$ws.Cells.Item(77, 8).Value = 2421.5715
$ws.Cells.Item(77, 9).Value = 1738.75
$ws.Cells.Item(77, 11).Value = 8693.75
$ws.Cells.Item(77, 13).Value = -4325.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 5424.9
$ws.Cells.Item(86, 9).Value = 1593.625
$ws.Cells.Item(86, 10).Value = 20750
$ws.Cells.Item(86, 11).Value = 1593.625
$ws.Cells.Item(86, 12).Value = 20750
$ws.Cells.Item(86, 13).Value = -470.625
$ws.Cells.Item(86, 14).Value = -22996
$ws.Cells.Item(89, 8).Value = 5424.9
$ws.Cells.Item(89, 9).Value = 1593.625
$ws.Cells.Item(89, 10).Value = 20750
$ws.Cells.Item(89, 11).Value = 7968.125
$ws.Cells.Item(89, 12).Value = 103750
$ws.Cells.Item(89, 13).Value = -2352.125
$ws.Cells.Item(89, 14).Value = -114982

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(8, 8).Value = 6500
$ws.Cells.Item(8, 10).Value = 10000
$ws.Cells.Item(8, 12).Value = 10000
$ws.Cells.Item(8, 14).Value = -10280
$ws.Cells.Item(31, 8).Value = 2020.1111
$ws.Cells.Item(31, 9).Value = 1599.2858
$ws.Cells.Item(31, 11).Value = 1599.2858
$ws.Cells.Item(31, 13).Value = -1304.2858
$ws.Cells.Item(34, 8).Value = 2020.1111
$ws.Cells.Item(34, 9).Value = 1599.2858
$ws.Cells.Item(34, 11).Value = 1599.2858
$ws.Cells.Item(34, 13).Value = -1397.2858
$ws.Cells.Item(43, 8).Value = 26367.5
$ws.Cells.Item(43, 10).Value = 26367.5
$ws.Cells.Item(43, 12).Value = 26367.5
$ws.Cells.Item(43, 14).Value = -26735.5
$ws.Cells.Item(101, 8).Value = 26367.5
$ws.Cells.Item(101, 10).Value = 26367.5
$ws.Cells.Item(101, 12).Value = 26367.5
$ws.Cells.Item(101, 14).Value = -32857.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(62, 8).Value = 6900
$ws.Cells.Item(62, 9).Value = 6900
$ws.Cells.Item(62, 11).Value = 20700
$ws.Cells.Item(62, 13).Value = -20014
$ws.Cells.Item(65, 8).Value = 6900
$ws.Cells.Item(65, 9).Value = 6900
$ws.Cells.Item(65, 11).Value = 62100
$ws.Cells.Item(65, 13).Value = -58668
$ws.Cells.Item(88, 8).Value = 18999
$ws.Cells.Item(88, 10).Value = 18999
$ws.Cells.Item(88, 12).Value = 56997
$ws.Cells.Item(88, 14).Value = -57853
$ws.Cells.Item(91, 8).Value = 18999
$ws.Cells.Item(91, 10).Value = 18999
$ws.Cells.Item(91, 12).Value = 56997
$ws.Cells.Item(91, 14).Value = -59961
$ws.Cells.Item(94, 8).Value = 2520.25
$ws.Cells.Item(94, 10).Value = 3027
$ws.Cells.Item(94, 12).Value = 9081
$ws.Cells.Item(94, 14).Value = -10433
$ws.Cells.Item(98, 8).Value = 422.5
$ws.Cells.Item(98, 10).Value = 545
$ws.Cells.Item(98, 12).Value = 1635
$ws.Cells.Item(98, 14).Value = -4631
$ws.Cells.Item(104, 8).Value = 25000
$ws.Cells.Item(104, 10).Value = 25000
$ws.Cells.Item(104, 12).Value = 75000
$ws.Cells.Item(104, 14).Value = -80242
$ws.Cells.Item(124, 8).Value = 7499.5
$ws.Cells.Item(124, 9).Value = 4999
$ws.Cells.Item(124, 11).Value = 14997
$ws.Cells.Item(124, 13).Value = -10087

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(18, 8).Value = 3500
$ws.Cells.Item(18, 9).Value = 3500
$ws.Cells.Item(18, 11).Value = 3500
$ws.Cells.Item(18, 13).Value = -3207
$ws.Cells.Item(58, 8).Value = 37250
$ws.Cells.Item(58, 9).Value = 37500
$ws.Cells.Item(58, 11).Value = 37500
$ws.Cells.Item(58, 13).Value = -37223
$ws.Cells.Item(132, 8).Value = 5386.375
$ws.Cells.Item(132, 9).Value = 5018.8
$ws.Cells.Item(132, 10).Value = 5999
$ws.Cells.Item(132, 11).Value = 15056.4
$ws.Cells.Item(132, 12).Value = 17997
$ws.Cells.Item(132, 13).Value = -12526.4
$ws.Cells.Item(132, 14).Value = -23057

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 619
$ws.Cells.Item(16, 9).Value = 619
$ws.Cells.Item(16, 11).Value = 619
$ws.Cells.Item(16, 13).Value = -449
$ws.Cells.Item(46, 8).Value = 3024.875
$ws.Cells.Item(46, 9).Value = 2549.75
$ws.Cells.Item(46, 11).Value = 2549.75
$ws.Cells.Item(46, 13).Value = -2361.75
$ws.Cells.Item(61, 8).Value = 1367.5834
$ws.Cells.Item(61, 9).Value = 1257.7778
$ws.Cells.Item(61, 10).Value = 1697
$ws.Cells.Item(61, 11).Value = 1257.7778
$ws.Cells.Item(61, 12).Value = 1697
$ws.Cells.Item(61, 13).Value = -1055.7778
$ws.Cells.Item(61, 14).Value = -2101
$ws.Cells.Item(100, 8).Value = 2099.2
$ws.Cells.Item(100, 9).Value = 2099.2
$ws.Cells.Item(100, 11).Value = 2099.2
$ws.Cells.Item(100, 13).Value = -1558.2
$ws.Cells.Item(113, 8).Value = 1367.5834
$ws.Cells.Item(113, 9).Value = 1257.7778
$ws.Cells.Item(113, 10).Value = 1697
$ws.Cells.Item(113, 11).Value = 1257.7778
$ws.Cells.Item(113, 12).Value = 1697
$ws.Cells.Item(113, 13).Value = 912.2221999999999
$ws.Cells.Item(113, 14).Value = -6037

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(32, 8).Value = 4005379.5
$ws.Cells.Item(32, 9).Value = 4005379.5
$ws.Cells.Item(32, 11).Value = 4005379.5
$ws.Cells.Item(32, 13).Value = -4005062.5
$ws.Cells.Item(62, 8).Value = 333335680
$ws.Cells.Item(62, 9).Value = 500002000
$ws.Cells.Item(62, 10).Value = 3000
$ws.Cells.Item(62, 11).Value = 500002000
$ws.Cells.Item(62, 12).Value = 3000
$ws.Cells.Item(62, 13).Value = -500001376
$ws.Cells.Item(62, 14).Value = -4248
$ws.Cells.Item(65, 8).Value = 333335680
$ws.Cells.Item(65, 9).Value = 500002000
$ws.Cells.Item(65, 10).Value = 3000
$ws.Cells.Item(65, 11).Value = 2500010000
$ws.Cells.Item(65, 12).Value = 15000
$ws.Cells.Item(65, 13).Value = -2500006880
$ws.Cells.Item(65, 14).Value = -21240
$ws.Cells.Item(105, 8).Value = 46765.668
$ws.Cells.Item(105, 10).Value = 46765.668
$ws.Cells.Item(105, 12).Value = 46765.668
$ws.Cells.Item(105, 14).Value = -53753.668
$ws.Cells.Item(122, 8).Value = 1087.25
$ws.Cells.Item(122, 9).Value = 1042.5714
$ws.Cells.Item(122, 11).Value = 3127.7142
$ws.Cells.Item(122, 13).Value = -677.7142000000003
$ws.Cells.Item(132, 8).Value = 2622.1052
$ws.Cells.Item(132, 9).Value = 1886.7858
$ws.Cells.Item(132, 11).Value = 5660.357400000001
$ws.Cells.Item(132, 13).Value = -3130.357400000001
